# Applies the "New crime data collected" update to the 112th Precinct
# weekly CompStat workbook:
#   - bumps the report Volume/Number and the covered date range
#   - refreshes the crime-complaint statistics table (rows 16-27)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text (rich-text shared strings) - edit only the affected runs
# ---------------------------------------------------------------------
# A8:  "Volume 30   Number  37" -> "...  38"
$ws.Range("A8").Characters(21, 2).Text = "38"

# C9: "Report Covering the Week  9/11/2023  Through  9/17/2023"
#  -> "Report Covering the Week  9/18/2023  Through  9/24/2023"
$ws.Range("C9").Characters(27, 9).Text = "9/18/2023"
$ws.Range("C9").Characters(47, 9).Text = "9/24/2023"

# ---------------------------------------------------------------------
# Helper: convert a cell that currently holds the shared placeholder
# text "0" or "***.*" into a real number, matching the number-format
# (and therefore cell style) of a known-good template cell.
# ---------------------------------------------------------------------
function Set-NumberFromText($cellAddr, $value, $templateAddr) {
    $cell = $ws.Range($cellAddr)
    $cell.NumberFormat = $ws.Range($templateAddr).NumberFormat
    $cell.Value = $value
}

# Helper: convert a numeric cell into the shared placeholder text
# ("0" or "***.*"), matching the style of a known-good template cell.
function Set-TextFromNumber($cellAddr, $text, $templateAddr) {
    $cell = $ws.Range($cellAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $ws.Range($templateAddr).Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 83.333333333333
$ws.Range("I16").Value = 54
$ws.Range("J16").Value = 67
$ws.Range("K16").Value = -19.402985074626
$ws.Range("L16").Value = 80
$ws.Range("M16").Value = -33.333333333333
$ws.Range("N16").Value = -88.486140724946

# ---------------------------------------------------------------------
# Row 17  (D17, E17 switch from text placeholders to real numbers)
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 4
Set-NumberFromText "D17" 2   "D16"
Set-NumberFromText "E17" 100 "H15"
$ws.Range("F17").Value = 7
$ws.Range("H17").Value = 75
$ws.Range("I17").Value = 74
$ws.Range("J17").Value = 67
$ws.Range("K17").Value = 10.447761194029
$ws.Range("L17").Value = 60.869565217391
$ws.Range("M17").Value = 80.487804878048
$ws.Range("N17").Value = -16.853932584269

# ---------------------------------------------------------------------
# Row 18
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 85
$ws.Range("J18").Value = 84
$ws.Range("K18").Value = 1.190476190476
$ws.Range("L18").Value = 28.787878787878
$ws.Range("M18").Value = -6.593406593406
$ws.Range("N18").Value = -91.567460317460

# ---------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 66.666666666666
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = -15.789473684210
$ws.Range("I19").Value = 352
$ws.Range("J19").Value = 378
$ws.Range("K19").Value = -6.878306878306
$ws.Range("L19").Value = 43.089430894308
$ws.Range("M19").Value = 27.536231884058
$ws.Range("N19").Value = -52.303523035230

# ---------------------------------------------------------------------
# Row 20
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 110
$ws.Range("J20").Value = 72
$ws.Range("K20").Value = 52.777777777777
$ws.Range("L20").Value = 139.130434782609
$ws.Range("M20").Value = 50.684931506849
$ws.Range("N20").Value = -95.652173913043

# ---------------------------------------------------------------------
# Row 21 (bold "TOTAL" row)
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 40
$ws.Range("F21").Value = 73
$ws.Range("G21").Value = 75
$ws.Range("H21").Value = -2.666666666666
$ws.Range("I21").Value = 678
$ws.Range("J21").Value = 679
$ws.Range("K21").Value = -0.147275405007
$ws.Range("L21").Value = 53.741496598639
$ws.Range("M21").Value = 20
$ws.Range("N21").Value = -86.006191950464

# ---------------------------------------------------------------------
# Row 22  (D22, E22 switch from text placeholders to real numbers)
# ---------------------------------------------------------------------
Set-NumberFromText "D22" 2    "D16"
Set-NumberFromText "E22" -100 "H15"
$ws.Range("J22").Value = 27
$ws.Range("K22").Value = -18.518518518518
$ws.Range("L22").Value = 144.444444444444
$ws.Range("M22").Value = 37.5

# ---------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = 24.137931034482
$ws.Range("F24").Value = 126
$ws.Range("G24").Value = 157
$ws.Range("H24").Value = -19.745222929936
$ws.Range("I24").Value = 1130
$ws.Range("J24").Value = 1328
$ws.Range("K24").Value = -14.909638554216
$ws.Range("L24").Value = 23.362445414847
$ws.Range("M24").Value = 55.647382920110

# ---------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 450
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 125
$ws.Range("I25").Value = 185
$ws.Range("J25").Value = 142
$ws.Range("K25").Value = 30.281690140845
$ws.Range("L25").Value = 31.205673758865
$ws.Range("M25").Value = 20.129870129870

# ---------------------------------------------------------------------
# Row 26  (C26 becomes the "0" text placeholder; D26, E26 become numbers)
# ---------------------------------------------------------------------
Set-TextFromNumber "C26" "0" "C22"
Set-NumberFromText "D26" 1    "D16"
Set-NumberFromText "E26" -100 "H15"
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 14
$ws.Range("K26").Value = -35.714285714285
$ws.Range("L26").Value = -30.769230769230

# ---------------------------------------------------------------------
# Row 27  (D27 becomes the "0" text placeholder; E27 becomes "***.*")
# ---------------------------------------------------------------------
Set-TextFromNumber "D27" "0"     "C22"
Set-TextFromNumber "E27" "***.*" "E23"
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -80

$excel.CutCopyMode = $false
